$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2 existing entry gets overwritten with real data
$ws.Range("B2").Value = "Zona 4"
$ws.Range("C2").Value = "CORONITAS SERVCOM SRL"
$ws.Range("D2").Value = 4
Set-TextValue $ws.Range("E2") "2025-08-11"

# Row 3: new entry
$ws.Range("A3").Value = "ioana"
$ws.Range("B3").Value = "Zona 3"
$ws.Range("C3").Value = "DEDEMAN"
$ws.Range("D3").Value = 7
Set-TextValue $ws.Range("E3") "2025-08-11"

# Row 4: new entry
$ws.Range("A4").Value = "andrei"
$ws.Range("B4").Value = "Zona 4"
$ws.Range("C4").Value = "AGROLIV SRL"
$ws.Range("D4").Value = 7
Set-TextValue $ws.Range("E4") "2025-08-11"

# Row 5: new entry
$ws.Range("A5").Value = "andrei"
$ws.Range("B5").Value = "Zona 5"
$ws.Range("C5").Value = "EUROCONSTRUCT SRL COVASNA"
$ws.Range("D5").Value = 100
Set-TextValue $ws.Range("E5") "2025-08-11"

# Row 6: new entry
$ws.Range("A6").Value = "ioana"
$ws.Range("B6").Value = "Zona 4"
$ws.Range("C6").Value = "MACON SRL BISTRITA"
$ws.Range("D6").Value = 677
Set-TextValue $ws.Range("E6") "2025-08-06"
